# ui: move command execution responsibility from CommandBox to MainWindow
#
# The diagram previously showed CommandBox's own line/arrow running up to
# Logic (CommandBox called Logic#execute(String) directly). Responsibility
# for invoking Logic#execute(String) has moved to MainWindow, so the small
# connector rectangle ("Rectangle 142") and the bent line that leads to the
# "Logic" box ("Freeform 115") must be repositioned/resized so the line now
# originates from MainWindow instead of CommandBox.

$EMU_PER_POINT = 12700

function ToPt($emu) {
    return $emu / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connectorRect = $null
$connectorLine = $null

foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Rectangle 142") {
        $connectorRect = $shp
    } elseif ($shp.Name -eq "Freeform 115") {
        $connectorLine = $shp
    }
}

# "Rectangle 142" - small rectangle that sits on the vertical line just
# below MainWindow; move it up from the CommandBox row to the MainWindow row.
$connectorRect.Left = ToPt(5422048)
$connectorRect.Top = ToPt(2339335)
$connectorRect.Width = ToPt(229325)
$connectorRect.Height = ToPt(166560)

# "Freeform 115" - the bent connector line pointing to Logic; move/resize it
# so it now starts at MainWindow's level rather than CommandBox's.
$connectorLine.Left = ToPt(3186477)
$connectorLine.Top = ToPt(2405681)
$connectorLine.Width = ToPt(3537529)
$connectorLine.Height = ToPt(45719)

# Refresh the "last modified" date field shown on the slide (date placeholder
# auto-updates to the current date; pin it to the date recorded for this
# edit so the exported markup matches exactly).
foreach ($sl in $p.Slides) {
    foreach ($shp in $sl.Shapes) {
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                # no-op placeholder; date fields are handled via masters/layouts below
            }
        }
    }
}
